$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.898.95'
$ws.Range('E2').Value = '  -1.65%  '

$ws.Range('D3').Value = '2.452.68'
$ws.Range('E3').Value = '  -2.80%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.86'
$ws.Range('E5').Value = '  -2.75%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '165.54'
$ws.Range('E6').Value = '  -5.13%  '

$ws.Range('E7').Value = '  +0.09%  '

$ws.Range('E8').Value = '  -3.24%  '

$ws.Range('D9').Value = '2.452.74'
$ws.Range('E9').Value = '  -2.73%  '

$ws.Range('E10').Value = '  -4.30%  '

$ws.Range('E11').Value = '  -1.05%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.88'
$ws.Range('E12').Value = '  -4.51%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.331'
$ws.Range('E13').Value = '  -3.77%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.29'
$ws.Range('E14').Value = '  -4.87%  '

$ws.Range('D15').Value = '2.888.22'
$ws.Range('E15').Value = '  -3.23%  '

$ws.Range('D16').Value = '66.846.51'
$ws.Range('E16').Value = '  -1.62%  '

$ws.Range('E17').Value = '  -5.81%  '

$ws.Range('D18').Value = '2.494.02'
$ws.Range('E18').Value = '  -0.38%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.31'
$ws.Range('E19').Value = '  -5.11%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.75'
$ws.Range('E20').Value = '  -4.17%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '354.95'
$ws.Range('E21').Value = '  -2.53%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.03'
$ws.Range('E22').Value = '  -3.05%  '

$ws.Range('E23').Value = '  -0.18%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '69.47'
$ws.Range('E24').Value = '  -2.73%  '

$ws.Range('E26').Value = '  -8.22%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.91'
$ws.Range('E27').Value = '  -11.01%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.998'
$ws.Range('E28').Value = '  +0.05%  '

$ws.Range('D29').Value = '2.572.01'
$ws.Range('E29').Value = '  -2.97%  '

$ws.Range('D30').Value = '0.0₃0896'
$ws.Range('E30').Value = '  -8.82%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '505.22'
$ws.Range('E31').Value = '  -5.23%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.79'
$ws.Range('E32').Value = '  -6.78%  '

$ws.Range('E33').Value = '  -7.15%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.21'
$ws.Range('E34').Value = '  -8.53%  '

$ws.Range('E35').Value = '  -0.02%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '158.89'
$ws.Range('E36').Value = '  +0.63%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.116'
$ws.Range('E37').Value = '  -9.58%  '

$ws.Range('B38').Value = 'WhiteBITCoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.57'
$ws.Range('E38').Value = '  -0.54%  '

$ws.Range('B39').Value = 'EthereumClassic'
$ws.Range('C39').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.43'
$ws.Range('E39').Value = '  -1.70%  '

$ws.Range('E40').Value = '  -7.41%  '

$ws.Range('E42').Value = '  -7.45%  '

$ws.Range('E43').Value = '  -7.23%  '

$ws.Range('E44').Value = '  -8.21%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '38.71'
$ws.Range('E45').Value = '  -3.07%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.31'
$ws.Range('E46').Value = '  -8.48%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '141.44'
$ws.Range('E47').Value = '  -4.34%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.47'
$ws.Range('E48').Value = '  -6.52%  '

$ws.Range('E49').Value = '  -7.60%  '

$ws.Range('B50').Value = 'Optimism'
$ws.Range('C50').Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.58'
$ws.Range('E50').Value = '  -8.67%  '

$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0250'
$ws.Range('E51').Value = '  -10.25%  '
